# Apply the "Add files via upload" edit to Sheet1:
#  - Extend the stock-screener table from 5 data rows (rows 2-6) to 8 data
#    rows (rows 2-9), refreshing the ticker lists in columns B..E.
#  - New rows 7-9 get the same index-column (A) styling as the existing
#    rows, obtained by copying the format from row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the index column (A) and replicate its style into the new rows ---
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7:A9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7

# Clear the rest of the new rows (B/C/E/F) so they read back blank, matching
# the shape of the existing rows' mostly-empty columns.
$ws.Range("B7:F9").Value = ""

# --- Row 2 ---
$ws.Range("B2").Value = "NSE:3MINDIA"
$ws.Range("C2").Value = "NSE:ADORWELD"
$ws.Range("D2").Value = "NSE:APOLLOTYRE"
$ws.Range("E2").Value = "NSE:AMBUJACEM"

# --- Row 3 ---
$ws.Range("B3").Value = "NSE:CARTRADE"
$ws.Range("C3").Value = "NSE:DCMSHRIRAM"
$ws.Range("D3").Value = "NSE:CHAMBLFERT"

# --- Row 4 ---
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "NSE:COLPAL"

# --- Row 5 ---
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "NSE:CONCOR"

# --- Row 6 ---
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "NSE:LICHSGFIN"

# --- Row 7 (new) ---
$ws.Range("D7").Value = "NSE:LTTS"

# --- Row 8 (new) ---
$ws.Range("D8").Value = "NSE:NESTLEIND"

# --- Row 9 (new) ---
$ws.Range("D9").Value = "NSE:ONGC"
